# Fruta / hortaliza, semanal
# Insert 5 new daily-price rows (American Nectar x2, Carson x3) ahead of the
# existing row 577 block, pushing the rest of the table down by 5 rows
# (dimension grows from A1:T609 to A1:T614).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows at 577..581 (existing rows 577-609 shift to 582-614).
$ws.Range("A577:A581").EntireRow.Insert()

# Common columns shared by every row in this block.
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$prodId    = 100103
$producto  = "Frutos de hueso (carozo)"
$catId     = 100103004
$categoria = "Durazno"

function Set-DuraznoRow($row, $fecha, $variedad, $calidad, $volumen, $pmin, $pmax, $pprom, $unidad, $origen, $pkg, $kgUnidad) {
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $prodId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $catId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $pmin
    $ws.Cells.Item($row, 15).Value = $pmax
    $ws.Cells.Item($row, 16).Value = $pprom
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $pkg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-DuraznoRow 577 44578 "American Nectar" "Especial" 330 18000 18000 18000 "`$/caja 18 kilos granel" "Región de O'Higgins" 1000 18
Set-DuraznoRow 578 44578 "American Nectar" "Primera"  350 14400 14400 14400 "`$/caja 18 kilos granel" "Región de O'Higgins" 800  18
Set-DuraznoRow 579 44578 "Carson"           "Especial" 380 27000 27000 27000 "`$/caja 18 kilos granel" "Región de O'Higgins" 1500 18
Set-DuraznoRow 580 44578 "Carson"           "Primera"  310 21600 21600 21600 "`$/caja 18 kilos granel" "Región de O'Higgins" 1200 18
Set-DuraznoRow 581 44578 "Carson"           "Segunda"  350 18000 18000 18000 "`$/caja 18 kilos granel" "Región de O'Higgins" 1000 18
